# Update automatico via Actualizar 02-19-2021 13-16-15
#
# The "Actualizar" refresh macro re-checks availability for every monitored
# service and stamps column D ("Ultimo"/Fecha) with the timestamp of the
# check. Each run keeps a rolling history of three 14-row batches: the
# newest batch lands in rows 2-15, the previous "newest" batch is pushed
# down to rows 16-29, and the one before that is pushed down to rows 30-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Oldest batch (rows 30-43) takes on the timestamp that used to belong to
# the middle batch (rows 16-29).
$tsOldest = 44246.51022087963
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $tsOldest
}

# Middle batch (rows 16-29) takes on the timestamp that used to belong to
# the newest batch (rows 2-15).
$tsMiddle = 44246.53152002315
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $tsMiddle
}

# Newest batch (rows 2-15) is stamped with the freshly captured timestamp
# from this update run.
$tsNewest = 44246.55280859994
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $tsNewest
}
